# Atualização de bases das ligas, do dia: 31-03-2024 às 20:29
#
# This script re-applies a set of corrections to the "Peru Liga 1" match
# odds table:
#   - several pairs of rows had their data (columns B..AC) swapped because
#     the two fixtures were recorded in the wrong order;
#   - one trio of rows had their data rotated for the same reason;
#   - the last few rows (263-266, i.e. sheet rows 265-268) had stale /
#     placeholder data which is replaced with the final match results, and
#     the two trailing placeholder rows (267-268) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $a = $ws.Range("B$rowA`:AC$rowA").Value2
    $b = $ws.Range("B$rowB`:AC$rowB").Value2
    $ws.Range("B$rowA`:AC$rowA").Value2 = $b
    $ws.Range("B$rowB`:AC$rowB").Value2 = $a
}

# ---------------------------------------------------------------------
# 1) Simple pairwise swaps (column A / row id stays put, B..AC swap)
# ---------------------------------------------------------------------
Swap-Rows 61 62
Swap-Rows 156 157
Swap-Rows 175 177
Swap-Rows 180 181
Swap-Rows 252 253

# ---------------------------------------------------------------------
# 2) Three-way rotation: 183 -> 184 -> 186 -> 183
# ---------------------------------------------------------------------
$row183 = $ws.Range("B183:AC183").Value2
$row184 = $ws.Range("B184:AC184").Value2
$row186 = $ws.Range("B186:AC186").Value2

$ws.Range("B184:AC184").Value2 = $row183
$ws.Range("B186:AC186").Value2 = $row184
$ws.Range("B183:AC183").Value2 = $row186

# ---------------------------------------------------------------------
# 3) Rows 265 & 266: replace with the final (non-placeholder) match data
# ---------------------------------------------------------------------
$ws.Cells.Item(265, 2).Value = 7971188
$ws.Cells.Item(265, 5).Value = 45380.875
$ws.Cells.Item(265, 6).Value = "Cusco FC"
$ws.Cells.Item(265, 7).Value = "Union Comercio"
$ws.Cells.Item(265, 8).Value = 1
$ws.Cells.Item(265, 9).Value = 0
$ws.Cells.Item(265, 10).Value = "H"
$ws.Cells.Item(265, 11).Value = 1.5
$ws.Cells.Item(265, 12).Value = 4
$ws.Cells.Item(265, 13).Value = 7
$ws.Cells.Item(265, 14).Value = 1.363
$ws.Cells.Item(265, 15).Value = 4.5
$ws.Cells.Item(265, 16).Value = 10
$ws.Cells.Item(265, 17).Value = -1.5
$ws.Cells.Item(265, 18).Value = 2
$ws.Cells.Item(265, 19).Value = 1.85
$ws.Cells.Item(265, 20).Value = 2.75
$ws.Cells.Item(265, 21).Value = 1.925
$ws.Cells.Item(265, 22).Value = 1.925
$ws.Cells.Item(265, 23).Value = 0.363
$ws.Cells.Item(265, 24).Value = -1
$ws.Cells.Item(265, 25).Value = -1
$ws.Cells.Item(265, 26).Value = -1
$ws.Cells.Item(265, 27).Value = 0.8500000000000001
$ws.Cells.Item(265, 28).Value = -1
$ws.Cells.Item(265, 29).Value = 0.925

$ws.Cells.Item(266, 2).Value = 7971189
$ws.Cells.Item(266, 5).Value = 45381.67708333334
$ws.Cells.Item(266, 6).Value = "Atletico Grau"
$ws.Cells.Item(266, 7).Value = "Carlos Manucci"
$ws.Cells.Item(266, 8).Value = 3
$ws.Cells.Item(266, 9).Value = 0
$ws.Cells.Item(266, 10).Value = "H"
$ws.Cells.Item(266, 11).Value = 1.615
$ws.Cells.Item(266, 12).Value = 4
$ws.Cells.Item(266, 13).Value = 5
$ws.Cells.Item(266, 14).Value = 1.65
$ws.Cells.Item(266, 15).Value = 3.75
$ws.Cells.Item(266, 16).Value = 5.75
$ws.Cells.Item(266, 17).Value = -1
$ws.Cells.Item(266, 18).Value = 2.025
$ws.Cells.Item(266, 19).Value = 1.775
$ws.Cells.Item(266, 20).Value = 2.5
$ws.Cells.Item(266, 21).Value = 1.875
$ws.Cells.Item(266, 22).Value = 1.925
$ws.Cells.Item(266, 23).Value = 0.6499999999999999
$ws.Cells.Item(266, 24).Value = -1
$ws.Cells.Item(266, 25).Value = -1
$ws.Cells.Item(266, 26).Value = 1.025
$ws.Cells.Item(266, 27).Value = -1
$ws.Cells.Item(266, 28).Value = 0.875
$ws.Cells.Item(266, 29).Value = -1

# ---------------------------------------------------------------------
# 4) Drop the two now-obsolete trailing placeholder rows (267 & 268)
# ---------------------------------------------------------------------
$ws.Range("A267:A268").EntireRow.Delete()
